$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(34, 8).Value = 0
$ws.Cells.Item(34, 9).Value = 0
$ws.Cells.Item(34, 10).Value = 0
$ws.Cells.Item(34, 11).Value = 0
$ws.Cells.Item(34, 12).Value = 0
$ws.Cells.Item(34, 13).ClearContents()
$ws.Cells.Item(36, 8).Value = 0
$ws.Cells.Item(36, 9).Value = 0
$ws.Cells.Item(36, 10).Value = 0
$ws.Cells.Item(36, 11).Value = 0
$ws.Cells.Item(36, 12).Value = 0
$ws.Cells.Item(36, 13).ClearContents()
$ws.Cells.Item(88, 8).Value = 2799.8
$ws.Cells.Item(88, 9).Value = 3400
$ws.Cells.Item(88, 10).Value = 2649.75
$ws.Cells.Item(88, 11).Value = 3400
$ws.Cells.Item(88, 12).Value = 2649.75
$ws.Cells.Item(88, 13).Value = -2994
$ws.Cells.Item(88, 14).Value = -3461.75
$ws.Cells.Item(91, 8).Value = 2799.8
$ws.Cells.Item(91, 9).Value = 3400
$ws.Cells.Item(91, 10).Value = 2649.75
$ws.Cells.Item(91, 11).Value = 3400
$ws.Cells.Item(91, 12).Value = 2649.75
$ws.Cells.Item(91, 13).Value = -1996
$ws.Cells.Item(91, 14).Value = -5457.75
$ws.Cells.Item(93, 8).Value = 0
$ws.Cells.Item(93, 9).Value = 0
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 11).Value = 0
$ws.Cells.Item(93, 12).Value = 0
$ws.Cells.Item(93, 14).ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1477.5
$ws.Cells.Item(32, 9).Value = 1477.5
$ws.Cells.Item(32, 10).Value = 0
$ws.Cells.Item(32, 11).Value = 1477.5
$ws.Cells.Item(32, 12).Value = 0
$ws.Cells.Item(32, 13).Value = -1190.5
$ws.Cells.Item(45, 8).Value = 3225.25
$ws.Cells.Item(45, 9).Value = 1912
$ws.Cells.Item(45, 10).Value = 3663
$ws.Cells.Item(45, 11).Value = 1912
$ws.Cells.Item(45, 12).Value = 3663
$ws.Cells.Item(45, 13).Value = -1535
$ws.Cells.Item(45, 14).Value = -4417
$ws.Cells.Item(63, 8).Value = 6815.9
$ws.Cells.Item(63, 9).Value = 6660.5
$ws.Cells.Item(63, 10).Value = 7049
$ws.Cells.Item(63, 11).Value = 6660.5
$ws.Cells.Item(63, 12).Value = 7049
$ws.Cells.Item(63, 13).Value = -5974.5
$ws.Cells.Item(63, 14).Value = -8421
$ws.Cells.Item(66, 8).Value = 6815.9
$ws.Cells.Item(66, 9).Value = 6660.5
$ws.Cells.Item(66, 10).Value = 7049
$ws.Cells.Item(66, 11).Value = 33302.5
$ws.Cells.Item(66, 12).Value = 35245
$ws.Cells.Item(66, 13).Value = -29870.5
$ws.Cells.Item(66, 14).Value = -42109
$ws.Cells.Item(122, 8).Value = 3308
$ws.Cells.Item(122, 9).Value = 3308
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 9924
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -7474
$ws.Cells.Item(122, 14).ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(6, 8).Value = 7000000.5
$ws.Cells.Item(6, 9).Value = 7000000.5
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = 7000000.5
$ws.Cells.Item(6, 12).Value = 0
$ws.Cells.Item(6, 13).Value = -6999887.5
$ws.Cells.Item(7, 8).Value = 170.38776
$ws.Cells.Item(7, 9).Value = 201.18182
$ws.Cells.Item(7, 10).Value = 161.47368
$ws.Cells.Item(7, 11).Value = 201.18182
$ws.Cells.Item(7, 12).Value = 161.47368
$ws.Cells.Item(7, 13).Value = -88.18181999999999
$ws.Cells.Item(7, 14).Value = -387.47368
$ws.Cells.Item(22, 8).Value = 253
$ws.Cells.Item(22, 9).Value = 260
$ws.Cells.Item(22, 10).Value = 190
$ws.Cells.Item(22, 11).Value = 260
$ws.Cells.Item(22, 12).Value = 190
$ws.Cells.Item(22, 13).Value = 90
$ws.Cells.Item(22, 14).Value = -890
$ws.Cells.Item(62, 8).Value = 3370
$ws.Cells.Item(62, 9).Value = 2990
$ws.Cells.Item(62, 10).Value = 3750
$ws.Cells.Item(62, 11).Value = 2990
$ws.Cells.Item(62, 12).Value = 3750
$ws.Cells.Item(62, 13).Value = -2366
$ws.Cells.Item(62, 14).Value = -4998
$ws.Cells.Item(65, 8).Value = 3370
$ws.Cells.Item(65, 9).Value = 2990
$ws.Cells.Item(65, 10).Value = 3750
$ws.Cells.Item(65, 11).Value = 14950
$ws.Cells.Item(65, 12).Value = 18750
$ws.Cells.Item(65, 13).Value = -11830
$ws.Cells.Item(65, 14).Value = -24990
$ws.Cells.Item(86, 8).Value = 9497.4
$ws.Cells.Item(86, 9).Value = 8747
$ws.Cells.Item(86, 10).Value = 12499
$ws.Cells.Item(86, 11).Value = 8747
$ws.Cells.Item(86, 12).Value = 12499
$ws.Cells.Item(86, 13).Value = -7624
$ws.Cells.Item(88, 8).Value = 18855.428
$ws.Cells.Item(88, 9).Value = 0
$ws.Cells.Item(88, 10).Value = 18855.428
$ws.Cells.Item(88, 11).Value = 0
$ws.Cells.Item(88, 12).Value = 18855.428
$ws.Cells.Item(88, 14).Value = -19667.428
$ws.Cells.Item(89, 8).Value = 9497.4
$ws.Cells.Item(89, 9).Value = 8747
$ws.Cells.Item(89, 10).Value = 12499
$ws.Cells.Item(89, 11).Value = 43735
$ws.Cells.Item(89, 12).Value = 62495
$ws.Cells.Item(89, 13).Value = -38119
$ws.Cells.Item(91, 8).Value = 18855.428
$ws.Cells.Item(91, 9).Value = 0
$ws.Cells.Item(91, 10).Value = 18855.428
$ws.Cells.Item(91, 11).Value = 0
$ws.Cells.Item(91, 12).Value = 18855.428
$ws.Cells.Item(91, 14).Value = -21663.428
$ws.Cells.Item(92, 8).Value = 19999
$ws.Cells.Item(92, 9).Value = 0
$ws.Cells.Item(92, 10).Value = 19999
$ws.Cells.Item(92, 11).Value = 0
$ws.Cells.Item(92, 12).Value = 19999
$ws.Cells.Item(92, 14).Value = -24991
$ws.Cells.Item(132, 8).Value = 0
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).ClearContents()
$ws.Cells.Item(132, 14).ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(32, 8).Value = 0
$ws.Cells.Item(32, 9).Value = 0
$ws.Cells.Item(32, 10).Value = 0
$ws.Cells.Item(32, 11).Value = 0
$ws.Cells.Item(32, 12).Value = 0
$ws.Cells.Item(32, 13).ClearContents()
$ws.Cells.Item(32, 14).ClearContents()
$ws.Cells.Item(41, 8).Value = 600
$ws.Cells.Item(41, 9).Value = 600
$ws.Cells.Item(41, 10).Value = 0
$ws.Cells.Item(41, 11).Value = 1800
$ws.Cells.Item(41, 12).Value = 0
$ws.Cells.Item(41, 13).Value = -1462
$ws.Cells.Item(44, 8).Value = 0
$ws.Cells.Item(44, 9).Value = 0
$ws.Cells.Item(44, 10).Value = 0
$ws.Cells.Item(44, 11).Value = 0
$ws.Cells.Item(44, 12).Value = 0
$ws.Cells.Item(44, 13).ClearContents()
$ws.Cells.Item(44, 14).ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(101, 8).Value = 20000
$ws.Cells.Item(101, 9).Value = 0
$ws.Cells.Item(101, 10).Value = 20000
$ws.Cells.Item(101, 11).Value = 0
$ws.Cells.Item(101, 12).Value = 20000
$ws.Cells.Item(101, 14).Value = -26490
$ws.Cells.Item(104, 8).Value = 0
$ws.Cells.Item(104, 9).Value = 0
$ws.Cells.Item(104, 10).Value = 0
$ws.Cells.Item(104, 11).Value = 0
$ws.Cells.Item(104, 12).Value = 0
$ws.Cells.Item(104, 14).ClearContents()
$ws.Cells.Item(122, 8).Value = 700
$ws.Cells.Item(122, 9).Value = 700
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 2100
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = 350
$ws.Cells.Item(136, 8).Value = 22749.5
$ws.Cells.Item(136, 9).Value = 0
$ws.Cells.Item(136, 10).Value = 22749.5
$ws.Cells.Item(136, 11).Value = 0
$ws.Cells.Item(136, 12).Value = 68248.5
$ws.Cells.Item(136, 14).Value = -73348.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(9, 8).Value = 332.5
$ws.Cells.Item(9, 9).Value = 332.5
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 11).Value = 332.5
$ws.Cells.Item(9, 12).Value = 0
$ws.Cells.Item(9, 13).Value = -108.5
$ws.Cells.Item(14, 8).Value = 6100
$ws.Cells.Item(14, 9).Value = 5000
$ws.Cells.Item(14, 10).Value = 6650
$ws.Cells.Item(14, 11).Value = 5000
$ws.Cells.Item(14, 12).Value = 6650
$ws.Cells.Item(14, 13).Value = -4828
$ws.Cells.Item(14, 14).Value = -6994
$ws.Cells.Item(22, 8).Value = 3462.3333
$ws.Cells.Item(22, 9).Value = 1772.4
$ws.Cells.Item(22, 10).Value = 5574.75
$ws.Cells.Item(22, 11).Value = 1772.4
$ws.Cells.Item(22, 12).Value = 5574.75
$ws.Cells.Item(22, 13).Value = -1477.4
$ws.Cells.Item(22, 14).Value = -6164.75
$ws.Cells.Item(27, 8).Value = 3462.3333
$ws.Cells.Item(27, 9).Value = 1772.4
$ws.Cells.Item(27, 10).Value = 5574.75
$ws.Cells.Item(27, 11).Value = 1772.4
$ws.Cells.Item(27, 12).Value = 5574.75
$ws.Cells.Item(27, 13).Value = -1665.4
$ws.Cells.Item(27, 14).Value = -5788.75
$ws.Cells.Item(46, 8).Value = 2767.8635
$ws.Cells.Item(46, 9).Value = 2524.375
$ws.Cells.Item(46, 10).Value = 2907
$ws.Cells.Item(46, 11).Value = 2524.375
$ws.Cells.Item(46, 12).Value = 2907
$ws.Cells.Item(46, 13).Value = -2336.375
$ws.Cells.Item(46, 14).Value = -3283
$ws.Cells.Item(53, 8).Value = 2867000
$ws.Cells.Item(53, 9).Value = 2867000
$ws.Cells.Item(53, 10).Value = 0
$ws.Cells.Item(53, 11).Value = 2867000
$ws.Cells.Item(53, 12).Value = 0
$ws.Cells.Item(53, 13).Value = -2866482
$ws.Cells.Item(74, 8).Value = 25000
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 10).Value = 25000
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 12).Value = 25000
$ws.Cells.Item(74, 14).Value = -26996
$ws.Cells.Item(77, 8).Value = 25000
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 10).Value = 25000
$ws.Cells.Item(77, 11).Value = 0
$ws.Cells.Item(77, 12).Value = 75000
$ws.Cells.Item(77, 14).Value = -84984

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(105, 8).Value = 28945
$ws.Cells.Item(105, 9).Value = 0
$ws.Cells.Item(105, 10).Value = 28945
$ws.Cells.Item(105, 11).Value = 0
$ws.Cells.Item(105, 12).Value = 28945
$ws.Cells.Item(105, 14).Value = -35933
